$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L2").Value = "[0.25797550888140935, 0.3725173484261755]"
$ws.Range("M2").Value = 0.0000000001728872600637033
$ws.Range("N2").Value = 0.0000000001728872600637033
$ws.Range("P2").Value = "[-0.855368570346771, -0.5031579825569228]"
$ws.Range("Q2").Value = 0.00000007890327147208609
$ws.Range("R2").Value = 0.00000007890327147208609
$ws.Range("T2").Value = "[0.4440753353504542, 0.503283453441965]"
$ws.Range("X2").Value = 2.013213213213272
$ws.Range("Y2").Value = 3.422462462462573
